$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds a date serial that was bumped from 45189 (2023-09-20)
# to 45190 (2023-09-21) for every data row (rows 2 through 119).
for ($row = 2; $row -le 119; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    $current = $cell.Value2
    if ($current -eq 45189) {
        $cell.Value = 45190
    }
}
